$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab07")

# --- Data refresh: updated ILOSTAT figures for columns F,G,H (informal
# employment %, total/female/male) and N,O,P (25-54 labour-force status)
# across the affected country/region rows. ---
$updates = @{
    "F6" = 89.141000000000005
    "G6" = 89.9
    "H6" = 88.548000000000002
    "N6" = 65.998649999999998
    "O6" = 15.199154
    "P6" = 18.802070000000001
    "F10" = 34.901000000000003
    "G10" = 34.648000000000003
    "H10" = 35.106999999999999
    "N10" = 52.600037
    "O10" = 23.536591000000001
    "P10" = 23.863371999999998
    "F13" = 74.333556000000002
    "G13" = 76.146556000000004
    "H13" = 72.531778000000003
    "N13" = 67.038786000000002
    "O13" = 13.785012
    "P13" = 19.176203000000001
    "F31" = 84.269000000000005
    "G31" = 86.668000000000006
    "H31" = 82.171000000000006
    "N31" = 70.493742999999995
    "O31" = 7.9363859999999997
    "P31" = 21.569870999999999
    "F32" = 14.222
    "G32" = 5.84
    "H32" = 23.016999999999999
    "N32" = 83.971172999999993
    "O32" = 1.6521459999999999
    "P32" = 14.378831999999999
    "F38" = 75.432417000000001
    "G38" = 76.871082999999999
    "H38" = 74.846999999999994
    "N38" = 78.049736999999993
    "O38" = 3.1143209999999999
    "P38" = 18.835943
    "F47" = 93.756
    "G47" = 96.006
    "H47" = 91.558999999999997
    "N47" = 82.522205
    "O47" = 2.1879460000000002
    "P47" = 15.289866
    "F61" = 89.778923000000006
    "G61" = 92.056385000000006
    "H61" = 87.559230999999997
    "N61" = 82.919415999999998
    "O61" = 2.5107970000000002
    "P61" = 14.569784
    "F62" = 80.322948999999994
    "G62" = 81.687949000000003
    "H62" = 78.998333000000002
    "N62" = 74.190470000000005
    "O62" = 4.8674289999999996
    "P62" = 20.9421
    "F63" = 37.622138
    "G63" = 36.903137999999998
    "H63" = 38.029915000000003
    "N63" = 74.503746000000007
    "O63" = 2.7750689999999998
    "P63" = 22.728971999999999
    "F64" = 52.890680000000003
    "G64" = 50.63776
    "H64" = 54.6614
    "N64" = 76.343633999999994
    "O64" = 3.9160789999999999
    "P64" = 19.740286999999999
    "F65" = 71.843050000000005
    "G65" = 71.110399999999998
    "H65" = 71.763999999999996
    "N65" = 69.831230000000005
    "O65" = 2.0036
    "P65" = 28.165081000000001
    "F66" = 50.143428999999998
    "G66" = 50.035525999999997
    "H66" = 50.043210999999999
    "N66" = 74.455906999999996
    "O66" = 3.0945830000000001
    "P66" = 22.456108
    "F67" = 73.047646999999998
    "G67" = 73.670647000000002
    "H67" = 72.421646999999993
    "N67" = 71.180952000000005
    "O67" = 3.1186859999999998
    "P67" = 25.700364
    "F68" = 85.493789000000007
    "G68" = 87.206737000000004
    "H68" = 83.850789000000006
    "N68" = 74.210514000000003
    "O68" = 2.965757
    "P68" = 22.823726000000001
    "F69" = 91.260999999999996
    "G69" = 93.631600000000006
    "H69" = 89.056600000000003
    "N69" = 76.847256000000002
    "O69" = 2.422129
    "P69" = 20.730615
    "F70" = 92.534999999999997
    "G70" = 94.72775
    "H70" = 90.491749999999996
    "N70" = 76.413148000000007
    "O70" = 3.3883350000000001
    "P70" = 20.198516999999999
    "F71" = 89.778923000000006
    "G71" = 92.056385000000006
    "H71" = 87.559230999999997
    "N71" = 82.919415999999998
    "O71" = 2.5107970000000002
    "P71" = 14.569784
    "F73" = 72.824332999999996
    "G73" = 73.414867000000001
    "H73" = 72.197599999999994
    "N73" = 70.845065000000005
    "O73" = 9.3947889999999994
    "P73" = 19.760147
    "F76" = 71.643570999999994
    "G76" = 71.772285999999994
    "H76" = 71.655000000000001
    "N76" = 80.942971
    "O76" = 1.310473
    "P76" = 17.746555000000001
    "F77" = 54.801614999999998
    "G77" = 53.507154
    "H77" = 55.623384999999999
    "N77" = 76.992750000000001
    "O77" = 4.6691349999999998
    "P77" = 18.338115999999999
    "F78" = 2.535593
    "G78" = 2.8571110000000002
    "H78" = 2.2556669999999999
    "N78" = 82.303816999999995
    "O78" = 4.6963670000000004
    "P78" = 12.999814000000001
    "F79" = 9.4875939999999996
    "G79" = 10.129875
    "H79" = 9.0340939999999996
    "N79" = 79.953193999999996
    "O79" = 3.7024180000000002
    "P79" = 16.344387000000001
    "F81" = 52.131875000000001
    "G81" = 47.878749999999997
    "H81" = 53.596874999999997
    "N81" = 72.650583999999995
    "O81" = 2.9982359999999999
    "P81" = 24.351179999999999
    "F82" = 79.723113999999995
    "G82" = 81.023143000000005
    "H82" = 78.427970999999999
    "N82" = 74.102256999999994
    "O82" = 4.7606830000000002
    "P82" = 21.137060000000002
    "F83" = 36.272395000000003
    "G83" = 35.882151
    "H83" = 36.581826
    "N83" = 74.686588
    "O83" = 2.75305
    "P83" = 22.568918
    "F84" = 92.718062000000003
    "G84" = 94.816999999999993
    "H84" = 90.662499999999994
    "N84" = 78.403726000000006
    "O84" = 2.2129240000000001
    "P84" = 19.38335
    "F86" = 78.876881999999995
    "G86" = 80.922528999999997
    "H86" = 76.953881999999993
    "N86" = 74.629189999999994
    "O86" = 3.8314859999999999
    "P86" = 21.539321999999999
    "F87" = 74.943888999999999
    "G87" = 74.6905
    "H87" = 74.855277999999998
    "N87" = 69.686667
    "O87" = 1.9551970000000001
    "P87" = 28.358027
    "F88" = 52.191249999999997
    "G88" = 49.947249999999997
    "H88" = 53.88
    "N88" = 53.580717999999997
    "O88" = 22.978207999999999
    "P88" = 23.441074
    "F89" = 49.622613000000001
    "G89" = 47.557386999999999
    "H89" = 50.778516000000003
    "N89" = 73.166820000000001
    "O89" = 3.4292769999999999
    "P89" = 23.433306000000002
    "F90" = 11.794364
    "G90" = 11.468363999999999
    "H90" = 12.101864000000001
    "N90" = 83.264769000000001
    "O90" = 3.3935559999999998
    "P90" = 13.341673999999999
    "F91" = 90.291269
    "G91" = 92.975499999999997
    "H91" = 88.131923
    "N91" = 77.358963000000003
    "O91" = 2.9361790000000001
    "P91" = 19.704857000000001
    "F93" = 59.106400000000001
    "G93" = 57.348399999999998
    "H93" = 60.719000000000001
    "N93" = 74.189918000000006
    "O93" = 3.0185080000000002
    "P93" = 22.791574000000001
    "F94" = 47.785842000000002
    "G94" = 44.054420999999998
    "H94" = 50.597842
    "N94" = 77.660971000000004
    "O94" = 2.629051
    "P94" = 19.709983000000001
    "F95" = 86.646083000000004
    "G95" = 88.456333000000001
    "H95" = 85.074250000000006
    "N95" = 77.612714999999994
    "O95" = 3.1140119999999998
    "P95" = 19.273274000000001
    "F96" = 62.347273000000001
    "G96" = 62.116
    "H96" = 62.626908999999998
    "N96" = 70.418625000000006
    "O96" = 2.8206669999999998
    "P96" = 26.758692
    "F97" = 88.406295999999998
    "G97" = 91.305593000000002
    "H97" = 86.029518999999993
    "N97" = 79.867075
    "O97" = 3.0619320000000001
    "P97" = 17.070992
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- Fix mojibake (double UTF-8 decoded as Latin-1) in the regional
# economic communities footnote so accented characters render correctly. ---
$ws.Range("A103").Value = 'Regional Economic Communities: CEN-SAD = "Community of Sahel-Saharan States"; COMESA = "Common Market for Eastern and Southern Africa"; EAC = "East African Community"; ECCAS = "Economic Community of Central African States"; ECOWAS = "Economic Community of West African States"; IGAD = "Intergovernmental Authority on Development"; SADC = "Southern African Development Community"; UMA = "Arab Maghreb Union"; PALOP = "Países Africanos de Língua Oficial Portuguesa"; ASEAN = "Association of Southeast Asian Nations"; MERCOSUR = "Mercado Común del Sur". EU27 = "European Union (27 members)". OECD = "Organisation for Economic Co-operation and Development".'
